$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 6
$ws.Range("A6").Value = 2
$ws.Range("K6").Value = 4
$ws.Range("O6").Value = 4

# Update values in row 8
$ws.Range("A8").Value = 2
$ws.Range("K8").Value = 3
$ws.Range("O8").Value = 4

# Update the selected cell/range to match the edited workbook
$ws.Range("V10").Select()
